# Update "想去人数" (F column) counts across the "展览" (sheet1), "演出"
# (sheet2), and "全部类型" (sheet4) worksheets to match the refreshed data
# snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 787
$ws1.Range("F4").Value = 753
$ws1.Range("F6").Value = 421
$ws1.Range("F7").Value = 652
$ws1.Range("F9").Value = 1246
$ws1.Range("F10").Value = 670
$ws1.Range("F11").Value = 400
$ws1.Range("F15").Value = 730
$ws1.Range("F18").Value = 363
$ws1.Range("F22").Value = 609

# --- 演出 (Performances) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 233

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 787
$ws4.Range("F7").Value = 753
$ws4.Range("F9").Value = 421
$ws4.Range("F10").Value = 652
$ws4.Range("F12").Value = 1246
$ws4.Range("F13").Value = 670
$ws4.Range("F16").Value = 400
$ws4.Range("F21").Value = 730
$ws4.Range("F25").Value = 363
$ws4.Range("F27").Value = 233
$ws4.Range("F35").Value = 609
